$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move AUTO.ROLLOVER header from F1 to M1, and insert the new headers
# in between (F1:L1), plus two more trailing headers (N1:O1).
$ws.Range("F1").Value = "INTEND.DATE"
$ws.Range("G1").Value = "CUST.REMARKS:1"
$ws.Range("H1").Value = "TAX.INTEREST.TYPE:1"
$ws.Range("I1").Value = "DRAWDOWN.ACCOUNT"
$ws.Range("J1").Value = "PRIN.LIQ.ACCT"
$ws.Range("K1").Value = "INT.LIQ.ACCT"
$ws.Range("L1").Value = "CHRG.LIQ.ACCT"
$ws.Range("M1").Value = "AUTO.ROLLOVER"
$ws.Range("N1").Value = "FINAL.MATURITY"
$ws.Range("O1").Value = "EXP.DATE"

# Row 2: the AUTO.ROLLOVER value (2) moves from F2 to M2; F2 becomes empty.
$ws.Range("F2").ClearContents()
$ws.Range("M2").Value = 2

# Match the author's recorded selection after the edit.
$ws.Range("F6").Select()
